# Generate Report for Handback
# Update the generated/handback timestamps on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 11:09:33"

# zh-cn sheet: handoff + handback datetimes for the first file row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 11:09:28"
$wsZhCn.Range("K2").Value = "2016-09-03 11:09:46"

# de-de sheet: handoff + handback datetimes for the first file row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 11:09:33"
$wsDeDe.Range("K2").Value = "2016-09-03 11:09:53"
